# Auto-generated Excel COM-interop script
# (1) normalizes shared-string text across several sheets (strip stray
#     spaces / full-width commas / dashes introduced by earlier OCR-ish import)
# (2) appends new metadata columns I:O to the "land" (土地) sheet, matching
#     the layout already used on the "stock" (股票) sheet
$wb = $excel.ActiveWorkbook

# ---- sheet 1: 土地 ----
$ws = $wb.Worksheets.Item(1)
# new columns: clone style from an existing header/data cell, then set value
$ws.Range("B1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "property_category"
$ws.Range("B1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "category"
$ws.Range("B1").Copy($ws.Range("K1"))
$ws.Range("K1").Value = "date"
$ws.Range("B1").Copy($ws.Range("L1"))
$ws.Range("L1").Value = "legislator_name"
$ws.Range("B1").Copy($ws.Range("M1"))
$ws.Range("M1").Value = "legislator_id"
$ws.Range("B1").Copy($ws.Range("N1"))
$ws.Range("N1").Value = "source_file"
$ws.Range("B1").Copy($ws.Range("O1"))
$ws.Range("O1").Value = "index"
$ws.Range("H2").Copy($ws.Range("I2"))
$ws.Range("I2").Value = "land"
$ws.Range("H2").Copy($ws.Range("J2"))
$ws.Range("J2").Value = "normal"
$ws.Range("H2").Copy($ws.Range("K2"))
$ws.Range("K2").Value = "2013-12-26"
$ws.Range("H2").Copy($ws.Range("L2"))
$ws.Range("L2").Value = "丁守中"
$ws.Range("H2").Copy($ws.Range("M2"))
$ws.Range("M2").Value = 515
$ws.Range("H2").Copy($ws.Range("N2"))
$ws.Range("N2").Value = "tmpc7fb1"
$ws.Range("H2").Copy($ws.Range("O2"))
$ws.Range("O2").Value = 13
$ws.Range("H3").Copy($ws.Range("I3"))
$ws.Range("I3").Value = "land"
$ws.Range("H3").Copy($ws.Range("J3"))
$ws.Range("J3").Value = "normal"
$ws.Range("H3").Copy($ws.Range("K3"))
$ws.Range("K3").Value = "2013-12-26"
$ws.Range("H3").Copy($ws.Range("L3"))
$ws.Range("L3").Value = "丁守中"
$ws.Range("H3").Copy($ws.Range("M3"))
$ws.Range("M3").Value = 515
$ws.Range("H3").Copy($ws.Range("N3"))
$ws.Range("N3").Value = "tmpc7fb1"
$ws.Range("H3").Copy($ws.Range("O3"))
$ws.Range("O3").Value = 14
$ws.Range("H4").Copy($ws.Range("I4"))
$ws.Range("I4").Value = "land"
$ws.Range("H4").Copy($ws.Range("J4"))
$ws.Range("J4").Value = "normal"
$ws.Range("H4").Copy($ws.Range("K4"))
$ws.Range("K4").Value = "2013-12-26"
$ws.Range("H4").Copy($ws.Range("L4"))
$ws.Range("L4").Value = "丁守中"
$ws.Range("H4").Copy($ws.Range("M4"))
$ws.Range("M4").Value = 515
$ws.Range("H4").Copy($ws.Range("N4"))
$ws.Range("N4").Value = "tmpc7fb1"
$ws.Range("H4").Copy($ws.Range("O4"))
$ws.Range("O4").Value = 15
$ws.Range("H5").Copy($ws.Range("I5"))
$ws.Range("I5").Value = "land"
$ws.Range("H5").Copy($ws.Range("J5"))
$ws.Range("J5").Value = "normal"
$ws.Range("H5").Copy($ws.Range("K5"))
$ws.Range("K5").Value = "2013-12-26"
$ws.Range("H5").Copy($ws.Range("L5"))
$ws.Range("L5").Value = "丁守中"
$ws.Range("H5").Copy($ws.Range("M5"))
$ws.Range("M5").Value = 515
$ws.Range("H5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = "tmpc7fb1"
$ws.Range("H5").Copy($ws.Range("O5"))
$ws.Range("O5").Value = 16
# text normalization (style already correct, only fix the value)
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "area"
$ws.Range("D1").Value = "share_portion"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "acquire_value"
$ws.Range("B2").Value = "臺北市北投區振興段一小段00930000地號"
$ws.Range("D2").Value = "10000分之101"
$ws.Range("F2").Value = "78年10月11曰"
$ws.Range("B3").Value = "臺北市北投區振興段一小段00660000地號"
$ws.Range("D3").Value = "97300分之4170"
$ws.Range("F3").Value = "93年08月19曰"
$ws.Range("B4").Value = "臺北市北投區振興段一小段00930000地號"
$ws.Range("D4").Value = "10000分之84"
$ws.Range("F4").Value = "85年05月240"
$ws.Range("B5").Value = "臺北市北投區振興段一小段00210008地號"
$ws.Range("F5").Value = "97年07月22日"

# ---- sheet 2: 建物 ----
$ws = $wb.Worksheets.Item(2)
# text normalization (style already correct, only fix the value)
$ws.Range("B2").Value = "臺北市北投區振興段一小段12407000建號"
$ws.Range("F2").Value = "78年10月11曰"
$ws.Range("B3").Value = "臺北市北投區振興段一小段12496000建號"
$ws.Range("F3").Value = "78年10月11曰"
$ws.Range("B4").Value = "臺北市北投區振興段一小段12398000建號"
$ws.Range("D4").Value = "全部"
$ws.Range("F4").Value = "85年05月24日"
$ws.Range("B5").Value = "臺北市北投區振興段一小段12496000建號"
$ws.Range("F5").Value = "85年05月24日"
$ws.Range("B6").Value = "臺北市北投區振興段一小段12884000建號"
$ws.Range("F6").Value = "93年08月19曰"
$ws.Range("B7").Value = "臺北市北投區振興段一小段12891000建號"
$ws.Range("D7").Value = "100000分之4464"
$ws.Range("F7").Value = "93年08月19日"

# ---- sheet 3: 汽車 ----
$ws = $wb.Worksheets.Item(3)
# text normalization (style already correct, only fix the value)
$ws.Range("B2").Value = "LEXUSES350(客車）"
$ws.Range("E2").Value = "100年04月19曰"

# ---- sheet 4: 存款 ----
$ws = $wb.Worksheets.Item(4)
# text normalization (style already correct, only fix the value)
$ws.Range("G8").Value = "10000"

# ---- sheet 6: 具有相當價值之財產 ----
$ws = $wb.Worksheets.Item(6)
# text normalization (style already correct, only fix the value)
$ws.Range("C1").Value = "項件"
$ws.Range("B3").Value = "手錶珠寶"

# ---- sheet 7: 保險 ----
$ws = $wb.Worksheets.Item(7)
# text normalization (style already correct, only fix the value)
$ws.Range("E4").Value = "新20年期增值分紅年繳21252元"
$ws.Range("E5").Value = "新20年期增值分紅年繳18063元"

# ---- sheet 8: 事業投資 ----
$ws = $wb.Worksheets.Item(8)
# text normalization (style already correct, only fix the value)
$ws.Range("D2").Value = "臺北市北投路二段13號10樓之一11"
$ws.Range("F2").Value = "95年09月01日"

